# Rewrite the "RelayTable" sheet (Sheet2) from the old 3-column
# Net_A/Net_B/Relay layout into the new 6-column
# Relays/I/CO/CC/controlA/controlB layout, per the commit
# "able to write relay pin connction back to excel".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$wsHeaderSrc = $wb.Worksheets.Item(1)

# Final table contents (header + 14 data rows -> A1:F15)
$newData = @(
    @("Relays", "I",    "CO",        "CC",        "controlA", "controlB"),
    @("K0",     "K0I",  "FOVI_1",    "Float",     "cbit0",     "VCC"),
    @("K1",     "VIN",  "KVI",       "K0I",       "cbit1",     "VCC"),
    @("K2",     "K2I",  "R_DRAIN_A", "Float",     "cbit2",     "VCC"),
    @("K3",     "K3I",  "CL_A",      "FOVI_2_F",  "cbit3",     "VCC"),
    @("K4",     "DRAIN_F", "K3I",    "K2I",       "cbit4",     "VCC"),
    @("K5",     "K5I",  "R_DRAIN_A", "Float",     "cbit2",     "VCC"),
    @("K6",     "K6I",  "CL_A",      "FOVI_2_F",  "cbit3",     "VCC"),
    @("K7",     "K7I",  "FOVI_2_S",  "GND",       "cbit5",     "VCC"),
    @("K8",     "K8I",  "K6I",       "K5I",       "cbit6",     "VCC"),
    @("K9",     "DRAIN_S", "K7I",    "K8I",       "cbit7",     "VCC"),
    @("K10",    "K10I", "FOVI_3",    "Float",     "cbit8",     "VCC"),
    @("K11",    "ROVP", "R_ROVP_A",  "K10I",      "cbit9",     "VCC"),
    @("K12",    "K12I", "FOVI_4_F",  "Float",     "cbit10",    "VCC"),
    @("K13",    "CS_F", "R_CS_A",    "K12I",      "cbit3",     "VCC")
)

$oldRowCount = 18
$newRowCount = 15
$oldColCount = 3
$newColCount = 6

# Drop the rows/columns that won't exist anymore in the smaller/bigger
# table before writing the new values into place (old sheet was
# A1:C18, new one is A1:F15).
if ($oldRowCount -gt $newRowCount) {
    $extraRows = $ws.Range($ws.Cells.Item($newRowCount + 1, 1), $ws.Cells.Item($oldRowCount, $oldColCount))
    $extraRows.Delete()
}

# Write every cell of the new table.
for ($r = 0; $r -lt $newData.Count; $r++) {
    $rowValues = $newData[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# Match the header row styling used elsewhere in the workbook (bold,
# boxed, centered) by copying the format of Sheet1's header cell onto
# the new A1:F1 header row.
$wsHeaderSrc.Range("A1").Copy()
$ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $newColCount)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the view state as it was (RelayTable sheet selected/active).
$ws.Select() | Out-Null
$ws.Range("A1").Select() | Out-Null
